$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.122.51'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '2.384.66'
$ws.Range("E3").Value = '  +4.05%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.510'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.53%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.26'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.08%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.122'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '2.757.82'
$ws.Range("E15").Value = '  +4.16%  '
$ws.Range("D16").Value = '2.398.47'
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.814'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.51%  '
$ws.Range("D18").Value = '43.121.46'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.27%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0732'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.75%  '
$ws.Range("E34").Value = '  -1.52%  '
$ws.Range("E35").Value = '  +7.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("E38").Value = '  +2.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.94%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -34.87%  '
$ws.Range("D43").Value = '1.952.36'
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("E44").Value = '  +0.77%  '
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.84%  '
$ws.Range("D48").Value = '2.612.65'
$ws.Range("E48").Value = '  +4.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.56%  '
$ws.Range("E51").Value = '  +1.80%  '
